# Apply "Added articles for the first version" edit to the active sheet (List1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: clear the stray F11 value (it becomes an empty formula-less cell)
$ws.Range("F11").ClearContents()

# Row 13 - new entry (3.3./4.3. style date row), with date, From/To times, note, duration formula
$ws.Range("B13").Value = 41902
$ws.Range("C13").Value = 0.64236111111111105
$ws.Range("D13").Value = 0.68055555555555547
$ws.Range("E13").Value = "zinecker"
$ws.Range("F13").Formula = "=D13-C13"

# Row 14
$ws.Range("C14").Value = 0.78819444444444453
$ws.Range("D14").Value = 0.83333333333333337
$ws.Range("E14").Value = "zinecker - 3 tables WTF"
$ws.Range("F14").Formula = "=D14-C14"

# Row 15 (no duration formula in this row)
$ws.Range("C15").Value = 0.97916666666666663
$ws.Range("D15").Value = 0.03125
$ws.Range("E15").Value = "zinecker - tab 4-10"

# Row 16
$ws.Range("C16").Value = 0.65972222222222221
$ws.Range("D16").Value = 0.73611111111111116
$ws.Range("E16").Value = "zinecker - zbytek tab"

# Row 17
$ws.Range("C17").Value = 0.78125
$ws.Range("D17").Value = 0.80902777777777779
$ws.Range("E17").Value = "zinecker -text"

# Row 18 (no duration formula in this row)
$ws.Range("C18").Value = 0.94097222222222221
$ws.Range("E18").Value = "zinecker - text - kurziva…"

# Row 19
$ws.Range("B19").Value = 41904
$ws.Range("C19").Value = 0.90625
$ws.Range("D19").Value = 0.96875
$ws.Range("E19").Value = "úpravy celkově"

# Fill the F16:F19 duration formula as one operation (like dragging the fill
# handle down in Excel) so it is stored as a single shared formula, then
# clear F18 back out since row 18 has no "Do" time and no duration.
$ws.Range("F16:F19").Formula = "=D16-C16"
$ws.Range("F18").ClearContents()

# Move the active selection to B20, matching the recorded cursor position.
$ws.Range("B20").Select()
